$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-05-05T14:17:01+00:00"

# --- Sheet "Elements": update Binding Strength + Binding Value Set ---
$wsElem = $wb.Worksheets.Item("Elements")

# Author.role binding strength: required -> preferred
$wsElem.Range("X5").Value = "preferred"
# Author.specialty binding strength: required -> preferred
$wsElem.Range("X6").Value = "preferred"

# Author.specialty binding value set URL changed from DMP to CISIS value set
$wsElem.Range("Z6").Value = "https://mos.esante.gouv.fr/NOS/JDV_J01-XdsAuthorSpecialty-CISIS/FHIR/JDV-J01-XdsAuthorSpecialty-CISIS"

# Column Z (26) widens to fit the new, longer "Binding Value Set" URL
$wsElem.Columns.Item(26).ColumnWidth = 83
